$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string to a cell while forcing text storage,
# even for strings that otherwise look like pure numbers (e.g. '6.28').
# Uses a scratch cell far outside the used range, marks it as Text via
# NumberFormat, assigns the literal, then copies only the *value*
# (PasteSpecial xlPasteValues = -4163) onto the destination so the
# destination cell's own formatting/style is left untouched.
function Set-TextValue {
    param($range, [string]$text)
    $helper = $ws.Range("ZZ1000")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $helper.Clear() | Out-Null
}

$ws.Range("D2").Value = '58.738.01'
$ws.Range("E2").Value = '  -1.25%  '

$ws.Range("D3").Value = '2.633.87'
$ws.Range("E3").Value = '  -1.37%  '

$ws.Range("E4").Value = '  +0.15%  '

Set-TextValue $ws.Range("D5") '520.59'
$ws.Range("E5").Value = '  +0.56%  '

Set-TextValue $ws.Range("D6") '145.14'
$ws.Range("E6").Value = '  -3.47%  '

$ws.Range("E7").Value = '  +0.24%  '

Set-TextValue $ws.Range("D8") '0.576'
$ws.Range("E8").Value = '  -0.33%  '

$ws.Range("D9").Value = '2.644.28'
$ws.Range("E9").Value = '  -0.86%  '

Set-TextValue $ws.Range("D10") '6.28'
$ws.Range("E10").Value = '  -3.93%  '

$ws.Range("E11").Value = '  -2.57%  '

$ws.Range("E12").Value = '  -2.48%  '

$ws.Range("E13").Value = '  -0.45%  '

$ws.Range("D14").Value = '3.094.16'
$ws.Range("E14").Value = '  -0.03%  '

$ws.Range("D15").Value = '58.740.79'
$ws.Range("E15").Value = '  -0.99%  '

Set-TextValue $ws.Range("D16") '20.77'
$ws.Range("E16").Value = '  -3.35%  '

$ws.Range("E17").Value = '  -2.79%  '

$ws.Range("D18").Value = '2.635.82'
$ws.Range("E18").Value = '  -1.07%  '

Set-TextValue $ws.Range("D19") '348.32'
$ws.Range("E19").Value = '  -0.14%  '

Set-TextValue $ws.Range("D20") '4.45'
$ws.Range("E20").Value = '  -4.25%  '

Set-TextValue $ws.Range("D21") '10.21'
$ws.Range("E21").Value = '  -4.06%  '

Set-TextValue $ws.Range("D22") '6.14'
$ws.Range("E22").Value = '  -1.47%  '

$ws.Range("E23").Value = '  -0.01%  '

Set-TextValue $ws.Range("D24") '61.68'
$ws.Range("E24").Value = '  +0.33%  '

Set-TextValue $ws.Range("D25") '0.414'
$ws.Range("E25").Value = '  -3.49%  '

Set-TextValue $ws.Range("D26") '0.164'

Set-TextValue $ws.Range("D27") '0.997'
$ws.Range("E27").Value = '  +0.58%  '

$ws.Range("D28").Value = '0.0₃0804'
$ws.Range("E28").Value = '  -4.39%  '

Set-TextValue $ws.Range("D29") '7.00'
$ws.Range("E29").Value = '  -2.49%  '

Set-TextValue $ws.Range("D30") '0.998'
$ws.Range("E30").Value = '  +0.12%  '

Set-TextValue $ws.Range("D31") '6.24'
$ws.Range("E31").Value = '  -5.22%  '

$ws.Range("E32").Value = '  -1.72%  '

$ws.Range("E33").Value = '  -0.17%  '

Set-TextValue $ws.Range("D34") '149.01'
$ws.Range("E34").Value = '  -0.34%  '

Set-TextValue $ws.Range("D35") '0.988'
$ws.Range("E35").Value = '  -5.96%  '

$ws.Range("E36").Value = '  -3.18%  '

$ws.Range("E37").Value = '  -2.36%  '

Set-TextValue $ws.Range("D38") '36.62'
$ws.Range("E38").Value = '  +0.32%  '

Set-TextValue $ws.Range("D39") '0.839'
$ws.Range("E39").Value = '  -6.20%  '

$ws.Range("E40").Value = '  -2.76%  '

$ws.Range("E41").Value = '  -2.91%  '

Set-TextValue $ws.Range("D42") '279.81'
$ws.Range("E42").Value = '  -4.95%  '

$ws.Range("E43").Value = '  +0.01%  '

Set-TextValue $ws.Range("D44") '0.0984'
$ws.Range("E44").Value = '  -1.75%  '

Set-TextValue $ws.Range("D45") '0.601'
$ws.Range("E45").Value = '  -4.53%  '

Set-TextValue $ws.Range("D46") '19.57'
$ws.Range("E46").Value = '  -1.31%  '

$ws.Range("E47").Value = '  -4.90%  '

Set-TextValue $ws.Range("D48") '10.32'
$ws.Range("E48").Value = '  +0.55%  '

$ws.Range("E51").Value = '  -4.11%  '

# Row 49 / 50: Maker and VeChain swap rank position (with refreshed values).
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D49") '0.0229'
$ws.Range("E49").Value = '  -2.57%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.988.64'
$ws.Range("E50").Value = '  +0.19%  '

